# 2023-07-04 Tuesday, 15:35:22 Auto Push
# Appends a new check-in log row (row 2) to the attendance sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 303
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = "王*佑"
$ws.Range("D2").Value = "2023-07-04 15:27:23"
# E2 (out-time) is still blank for this check-in; touch a no-op format
# property so the cell is materialized in the sheet instead of being
# dropped as a fully-empty cell.
$ws.Range("E2").Font.Bold = $false
$ws.Range("F2").Value = "IN"
